$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price list values (Noviembre 2025)
$data = @{
    2  = @(92690, 119053, 156676, 265611, 458513, 80727, 69634)
    3  = @(93093, 127992, 158159, 272922, 459933, 80939, 69907)
    4  = @(104721, 128254, 178215, 307192, 517867, 90957, 78853)
    5  = @(147817, 211350, 247556, 414844, 617706, 128946, 112158)
    6  = @(254165, 338559, 445599, 621083, 846369, 222705, 194366)
    7  = @(66469, 91010, 118953, 224832, 328022, 63185, 56866)
    8  = @(67716, 91909, 118915, 226685, 334111, 62883, 58691)
    9  = @(77720, 104039, 137296, 263011, 385785, 73518, 68864)
    10 = @(94721, 133153, 163532, 316648, 451433, 88322, 82731)
    11 = @(78785, 101192, 133173, 225769, 389736, 68619, 59188)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}

# Update the selected range shown in the sheet view
$ws.Range("I19:J19").Select()
